$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 39) with the forecast values
$ws.Range("A39").Value = 45986
$ws.Range("B39").Value = 2025
$ws.Range("C39").Value = 0.2298740481777584
$ws.Range("D39").Value = 2026
$ws.Range("E39").Value = -0.05255865067609333

# Reuse the existing cell style from the row above (row 38) for the date cell
# in column A, so no new/duplicate style entries are introduced.
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
